$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 1.674957333333333
$ws.Range("H2").Value2 = 5.024872
$ws.Range("I2").Value2 = 0.007413057527797912
$ws.Range("J2").Value2 = 0.007413057527797913
$ws.Range("M2").Value2 = 10.718847
$ws.Range("N2").Value2 = 32.156541
$ws.Range("O2").Value2 = 0.2473266771098565
$ws.Range("P2").Value2 = 0.2473266771098565
$ws.Range("Q2").Value2 = 17.953611387528
$ws.Range("R2").Value2 = 161.582502487752
$ws.Range("S2").Value2 = 0.001833446885574465
$ws.Range("T2").Value2 = 0.001833446885574466
$ws.Range("G3").Value2 = 1.674957333333333
$ws.Range("H3").Value2 = 5.024872
$ws.Range("I3").Value2 = 0.007413057527797912
$ws.Range("J3").Value2 = 0.007413057527797913
$ws.Range("O3").Value2 = 0.2437114581515935
$ws.Range("P3").Value2 = 0.2437114581515935
$ws.Range("Q3").Value2 = 17.69118018917955
$ws.Range("R3").Value2 = 159.220621702616
$ws.Range("S3").Value2 = 0.001806647059461276
$ws.Range("T3").Value2 = 0.001806647059461276
$ws.Range("G4").Value2 = 1.674957333333333
$ws.Range("H4").Value2 = 5.024872
$ws.Range("I4").Value2 = 0.007413057527797912
$ws.Range("J4").Value2 = 0.007413057527797913
$ws.Range("M4").Value2 = 7.214691666666667
$ws.Range("N4").Value2 = 21.644075
$ws.Range("O4").Value2 = 0.1664717964804274
$ws.Range("P4").Value2 = 0.1664717964804274
$ws.Range("Q4").Value2 = 12.08430071482222
$ws.Range("R4").Value2 = 108.7587064334
$ws.Range("S4").Value2 = 0.001234065004065274
$ws.Range("T4").Value2 = 0.001234065004065274
$ws.Range("G5").Value2 = 1.674957333333333
$ws.Range("H5").Value2 = 5.024872
$ws.Range("I5").Value2 = 0.007413057527797912
$ws.Range("J5").Value2 = 0.007413057527797913
$ws.Range("M5").Value2 = 14.84311633333333
$ws.Range("N5").Value2 = 44.529349
$ws.Range("O5").Value2 = 0.3424900682581226
$ws.Range("P5").Value2 = 0.3424900682581225
$ws.Range("Q5").Value2 = 24.86158655203645
$ws.Range("R5").Value2 = 223.754278968328
$ws.Range("S5").Value2 = 0.002538898578696896
$ws.Range("T5").Value2 = 0.002538898578696896
$ws.Range("G6").Value2 = 5.848171333333333
$ws.Range("I6").Value2 = 0.02588294618833193
$ws.Range("J6").Value2 = 0.02588294618833193
$ws.Range("M6").Value2 = 10.718847
$ws.Range("N6").Value2 = 32.156541
$ws.Range("O6").Value2 = 0.2473266771098565
$ws.Range("P6").Value2 = 0.2473266771098565
$ws.Range("Q6").Value2 = 62.68565375178599
$ws.Range("R6").Value2 = 564.170883766074
$ws.Range("S6").Value2 = 0.006401543074573363
$ws.Range("T6").Value2 = 0.006401543074573364
$ws.Range("G7").Value2 = 5.848171333333333
$ws.Range("I7").Value2 = 0.02588294618833193
$ws.Range("J7").Value2 = 0.02588294618833193
$ws.Range("O7").Value2 = 0.2437114581515935
$ws.Range("P7").Value2 = 0.2437114581515935
$ws.Range("Q7").Value2 = 61.76936616606022
$ws.Range("R7").Value2 = 555.924295494542
$ws.Range("S7").Value2 = 0.006307970556817604
$ws.Range("T7").Value2 = 0.006307970556817605
$ws.Range("G8").Value2 = 5.848171333333333
$ws.Range("I8").Value2 = 0.02588294618833193
$ws.Range("J8").Value2 = 0.02588294618833193
$ws.Range("M8").Value2 = 7.214691666666667
$ws.Range("N8").Value2 = 21.644075
$ws.Range("O8").Value2 = 0.1664717964804274
$ws.Range("P8").Value2 = 0.1664717964804274
$ws.Range("Q8").Value2 = 42.19275298383889
$ws.Range("R8").Value2 = 379.73477685455
$ws.Range("S8").Value2 = 0.004308780550177847
$ws.Range("T8").Value2 = 0.004308780550177848
$ws.Range("G9").Value2 = 5.848171333333333
$ws.Range("I9").Value2 = 0.02588294618833193
$ws.Range("J9").Value2 = 0.02588294618833193
$ws.Range("M9").Value2 = 14.84311633333333
$ws.Range("N9").Value2 = 44.529349
$ws.Range("O9").Value2 = 0.3424900682581226
$ws.Range("P9").Value2 = 0.3424900682581225
$ws.Range("Q9").Value2 = 86.80508743793179
$ws.Range("R9").Value2 = 781.2457869413861
$ws.Range("S9").Value2 = 0.008864652006763117
$ws.Range("T9").Value2 = 0.008864652006763115
$ws.Range("G10").Value2 = 6.742607
$ws.Range("H10").Value2 = 20.227821
$ws.Range("I10").Value2 = 0.02984155630929478
$ws.Range("J10").Value2 = 0.02984155630929478
$ws.Range("M10").Value2 = 10.718847
$ws.Range("N10").Value2 = 32.156541
$ws.Range("O10").Value2 = 0.2473266771098565
$ws.Range("P10").Value2 = 0.2473266771098565
$ws.Range("Q10").Value2 = 72.27297281412899
$ws.Range("R10").Value2 = 650.4567553271609
$ws.Range("S10").Value2 = 0.007380612961764552
$ws.Range("T10").Value2 = 0.007380612961764552
$ws.Range("G11").Value2 = 6.742607
$ws.Range("H11").Value2 = 20.227821
$ws.Range("I11").Value2 = 0.02984155630929478
$ws.Range("J11").Value2 = 0.02984155630929478
$ws.Range("O11").Value2 = 0.2437114581515935
$ws.Range("P11").Value2 = 0.2437114581515935
$ws.Range("Q11").Value2 = 71.21654564444032
$ws.Range("R11").Value2 = 640.948910799963
$ws.Range("S11").Value2 = 0.007272729201651115
$ws.Range("T11").Value2 = 0.007272729201651115
$ws.Range("G12").Value2 = 6.742607
$ws.Range("H12").Value2 = 20.227821
$ws.Range("I12").Value2 = 0.02984155630929478
$ws.Range("J12").Value2 = 0.02984155630929478
$ws.Range("M12").Value2 = 7.214691666666667
$ws.Range("N12").Value2 = 21.644075
$ws.Range("O12").Value2 = 0.1664717964804274
$ws.Range("P12").Value2 = 0.1664717964804274
$ws.Range("Q12").Value2 = 48.64583053450833
$ws.Range("R12").Value2 = 437.812474810575
$ws.Range("S12").Value2 = 0.004967777488580134
$ws.Range("T12").Value2 = 0.004967777488580134
$ws.Range("G13").Value2 = 6.742607
$ws.Range("H13").Value2 = 20.227821
$ws.Range("I13").Value2 = 0.02984155630929478
$ws.Range("J13").Value2 = 0.02984155630929478
$ws.Range("M13").Value2 = 14.84311633333333
$ws.Range("N13").Value2 = 44.529349
$ws.Range("O13").Value2 = 0.3424900682581226
$ws.Range("P13").Value2 = 0.3424900682581225
$ws.Range("Q13").Value2 = 100.0813000909477
$ws.Range("R13").Value2 = 900.731700818529
$ws.Range("S13").Value2 = 0.01022043665729898
$ws.Range("T13").Value2 = 0.01022043665729898
$ws.Range("G14").Value2 = 211.6811596666666
$ws.Range("H14").Value2 = 635.0434789999999
$ws.Range("I14").Value2 = 0.9368624399745754
$ws.Range("J14").Value2 = 0.9368624399745754
$ws.Range("M14").Value2 = 10.718847
$ws.Range("N14").Value2 = 32.156541
$ws.Range("O14").Value2 = 0.2473266771098565
$ws.Range("P14").Value2 = 0.2473266771098565
$ws.Range("Q14").Value2 = 2268.977963249571
$ws.Range("R14").Value2 = 20420.80166924614
$ws.Range("S14").Value2 = 0.2317110741879441
$ws.Range("T14").Value2 = 0.2317110741879441
$ws.Range("G15").Value2 = 211.6811596666666
$ws.Range("H15").Value2 = 635.0434789999999
$ws.Range("I15").Value2 = 0.9368624399745754
$ws.Range("J15").Value2 = 0.9368624399745754
$ws.Range("O15").Value2 = 0.2437114581515935
$ws.Range("P15").Value2 = 0.2437114581515935
$ws.Range("Q15").Value2 = 2235.81190027377
$ws.Range("R15").Value2 = 20122.30710246393
$ws.Range("S15").Value2 = 0.2283241113336635
$ws.Range("T15").Value2 = 0.2283241113336635
$ws.Range("G16").Value2 = 211.6811596666666
$ws.Range("H16").Value2 = 635.0434789999999
$ws.Range("I16").Value2 = 0.9368624399745754
$ws.Range("J16").Value2 = 0.9368624399745754
$ws.Range("M16").Value2 = 7.214691666666667
$ws.Range("N16").Value2 = 21.644075
$ws.Range("O16").Value2 = 0.1664717964804274
$ws.Range("P16").Value2 = 0.1664717964804274
$ws.Range("Q16").Value2 = 1527.214298637436
$ws.Range("R16").Value2 = 13744.92868773692
$ws.Range("S16").Value2 = 0.1559611734376042
$ws.Range("T16").Value2 = 0.1559611734376042
$ws.Range("G17").Value2 = 211.6811596666666
$ws.Range("H17").Value2 = 635.0434789999999
$ws.Range("I17").Value2 = 0.9368624399745754
$ws.Range("J17").Value2 = 0.9368624399745754
$ws.Range("M17").Value2 = 14.84311633333333
$ws.Range("N17").Value2 = 44.529349
$ws.Range("O17").Value2 = 0.3424900682581226
$ws.Range("P17").Value2 = 0.3424900682581225
$ws.Range("Q17").Value2 = 3142.008078507241
$ws.Range("R17").Value2 = 28278.07270656517
$ws.Range("S17").Value2 = 0.3208660810153636
$ws.Range("T17").Value2 = 0.3208660810153635
